# Build site at 2022-09-26 16:07:08 UTC
#
# Rewrites the "8800008" discipline sheet (Projetos Especiais em Engenharia I):
#  - Row 10 (Objetivos:) now shows the responsible-professor line instead of
#    the long objectives paragraph.
#  - The old orphan row 13 (professor name, no label) is repurposed into
#    "Programa resumido:" / "Semestral".
#  - Short syllabus / Programa / Syllabus content is cleared out, and
#    Programa: now shows the activation date.
#  - Avaliacao: loses its custom row height.
#  - Metodo: gains the professor line; Criterio/Norma de
#    recuperacao/Bibliografia each shift up one slot.
#  - The old last row (22, full bibliography) is removed entirely so the
#    sheet ends at row 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 10: Objetivos:  (label unchanged) - swap the long paragraph for the
# professor line. B10/C10 already exist, so their style is untouched.
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = "198273 - Domingos Savio Giordani"
$ws.Cells.Item(10, 3).Value = "198273 - Domingos Savio Giordani"
$ws.Rows.Item(10).RowHeight = 60

# Row 11: Objectives:  - untouched
# Row 12: Docentes responsáveis:  - untouched

# ---------------------------------------------------------------------
# Row 13: used to be an orphan B/C-only row holding the professor name.
# Gains a new A13 label "Programa resumido:" and B13/C13 becomes "Semestral".
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------
# Row 14: becomes "Short syllabus:" label only - clear the long B/C text.
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "Short syllabus:"
$ws.Cells.Item(14, 2).Clear()
$ws.Cells.Item(14, 3).Clear()
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------
# Row 15: becomes "Programa:" / "01/01/2020". B15/C15 are brand-new cells,
# and "01/01/2020" looks like a date, so force text entry first, then fix
# the resulting style back to the plain column style via a format-only paste
# from row 14 (whose B/C already carry the correct, unmerged style).
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = "Programa:"

$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "01/01/2020"
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122)

$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "01/01/2020"
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)

$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------------
# Row 16: becomes "Syllabus:" label only - clear the long B/C text.
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = "Syllabus:"
$ws.Cells.Item(16, 2).Clear()
$ws.Cells.Item(16, 3).Clear()
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------
# Row 17: becomes "Avaliação:" label only, and loses its custom row height
# entirely (back to the sheet default).
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# ---------------------------------------------------------------------
# Row 18: becomes "Método:" / professor line. B18/C18 are brand-new cells;
# fix up B18's style (column B's <col> ranges overlap, so new cells default
# to the bold label style) via a format-only paste from row 19.
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 1).Value = "Método:"

$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(18, 2).PasteSpecial(-4122)
$ws.Cells.Item(18, 2).Value = "198273 - Domingos Savio Giordani"

$ws.Cells.Item(19, 3).Copy()
$ws.Cells.Item(18, 3).PasteSpecial(-4122)
$ws.Cells.Item(18, 3).Value = "198273 - Domingos Savio Giordani"

$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------------
# Row 19: becomes "Critério:" / "Seminários e Estudos de Casos."
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 1).Value = "Critério:"
$ws.Cells.Item(19, 2).Value = "Seminários e Estudos de Casos."
$ws.Cells.Item(19, 3).Value = "Seminários e Estudos de Casos."
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------
# Row 20: becomes "Norma de recuperação:" / the evaluation-board paragraph
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value = "Serão feitas duas avaliações por uma banca de professores que assistirão aos seminários apresentados, as notas serão as médias das notas dadas pelos professores."
$ws.Cells.Item(20, 3).Value = "Serão feitas duas avaliações por uma banca de professores que assistirão aos seminários apresentados, as notas serão as médias das notas dadas pelos professores."
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------
# Row 21: becomes "Bibliografia:" / the resubmission paragraph
# ---------------------------------------------------------------------
$ws.Cells.Item(21, 1).Value = "Bibliografia:"
$ws.Cells.Item(21, 2).Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."
$ws.Cells.Item(21, 3).Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------
# Row 22: old Bibliografia: / full bibliography text row is removed, the
# sheet now ends at row 21.
# ---------------------------------------------------------------------
$ws.Rows.Item(22).Delete()
